$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A ("Nº Pedido") holds order numbers that must stay as text, not numbers,
# just like they were originally (inline string cells). Force text format first
# so Excel doesn't silently convert the numeric-looking values to real numbers.
$ws.Range("A2:A17").NumberFormat = "@"

# Widen column E (Bairro) from 12 to 14 characters.
$ws.Columns.Item(5).ColumnWidth = 13.1667

# Row 2
$ws.Range("A2").Value = "11083"
$ws.Range("B2").Value = "19/11/2024"
$ws.Range("E2").Value = "Brejarú"

# Row 3
$ws.Range("A3").Value = "11083"
$ws.Range("B3").Value = "19/11/2024"
$ws.Range("E3").Value = "Brejarú"

# Row 4
$ws.Range("A4").Value = "11084"
$ws.Range("B4").Value = "19/11/2024"
$ws.Range("E4").Value = "Brejarú"

# Row 5
$ws.Range("A5").Value = "11084"
$ws.Range("B5").Value = "19/11/2024"
$ws.Range("E5").Value = "Brejarú"

# Row 6
$ws.Range("A6").Value = "11084"
$ws.Range("B6").Value = "19/11/2024"
$ws.Range("C6").Value = "manhã"
$ws.Range("E6").Value = "Brejarú"

# Row 7
$ws.Range("A7").Value = "11084"
$ws.Range("B7").Value = "19/11/2024"
$ws.Range("E7").Value = "Brejarú"

# Row 8
$ws.Range("A8").Value = "11104"
$ws.Range("B8").Value = "19/11/2024"
$ws.Range("E8").Value = "Caminho Novo"

# Row 9
$ws.Range("A9").Value = "11104"
$ws.Range("B9").Value = "19/11/2024"
$ws.Range("E9").Value = "Caminho Novo"

# Row 10
$ws.Range("A10").Value = "11104"
$ws.Range("B10").Value = "19/11/2024"
$ws.Range("C10").Value = "tarde"
$ws.Range("E10").Value = "Caminho Novo"

# Row 11
$ws.Range("A11").Value = "11104"
$ws.Range("B11").Value = "19/11/2024"
$ws.Range("C11").Value = "tarde"
$ws.Range("E11").Value = "Caminho Novo"

# Row 12
$ws.Range("A12").Value = "11184"
$ws.Range("B12").Value = "21/11/2024"
$ws.Range("E12").Value = "Pagani"

# Row 13
$ws.Range("A13").Value = "11184"
$ws.Range("B13").Value = "21/11/2024"
$ws.Range("E13").Value = "Pagani"

# Row 14
$ws.Range("A14").Value = "11184"
$ws.Range("B14").Value = "21/11/2024"
$ws.Range("C14").Value = "manhã"
$ws.Range("E14").Value = "Pagani"

# Row 15
$ws.Range("A15").Value = "11184"
$ws.Range("B15").Value = "21/11/2024"
$ws.Range("E15").Value = "Pagani"

# Row 16
$ws.Range("A16").Value = "11186"
$ws.Range("B16").Value = "21/11/2024"
$ws.Range("E16").Value = "Passa Vinte"

# Row 17
$ws.Range("A17").Value = "11186"
$ws.Range("B17").Value = "21/11/2024"
$ws.Range("E17").Value = "Passa Vinte"
